# Append the new "2021年" row (row 13) to Sheet1, mirroring the style and
# layout of the preceding "2020年" row (row 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

# Copy the formatting (bold font, border, centered alignment) of the year
# label cell in the row above so the new label picks up the same style.
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = "2021年"

$ws.Cells.Item($row, 2).Value = 25
$ws.Cells.Item($row, 3).Value = 228
$ws.Cells.Item($row, 4).Value = 5
$ws.Cells.Item($row, 5).Value = 211
$ws.Cells.Item($row, 6).Value = 32
$ws.Cells.Item($row, 7).Value = 110
$ws.Cells.Item($row, 8).Value = 3
$ws.Cells.Item($row, 9).Value = 1
# Column J (index 10) has no data for this year - leave it blank.
$ws.Cells.Item($row, 11).Value = 27
$ws.Cells.Item($row, 12).Value = 29
$ws.Cells.Item($row, 13).Value = 38
$ws.Cells.Item($row, 14).Value = 5
$ws.Cells.Item($row, 15).Value = 714
